$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Auftrag"
$ws.Range("B1").Value = "Kategorie"
$ws.Range("C1").Value = "Arbeitschritt"
$ws.Range("D1").Value = "Menge"
$ws.Range("E1").Value = "Maschiene"
$ws.Range("F1").Value = "Arbeitskraft"
$ws.Range("G1").Value = "Zeit"

$ws.Rows.Item(1).RowHeight = 42

$ws.Range("H1").Select()
